$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (pushes existing data rows down by one).
$ws.Rows.Item(2).Insert()

# The newly inserted row inherits the bold header formatting from row 1
# above it; strip that back off so row 2 matches the plain data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new data values
$ws.Range("A2").Value = -0.0018325957935303
$ws.Range("B2").Value = -0.0296269636601209
$ws.Range("C2").Value = -0.0087048299610614

# The insert pushed the former last two rows (old rows 21 and 22, now at rows 22 and 23)
# beyond the desired range, so remove them.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
